# Trade #140 closed at 2026-02-18 00:41:47 - unknown UNKNOWN +0.000%
#
# Helper to write a text value (dates / times / free text) into a cell
# while keeping it stored as plain text without picking up a date/time
# number format or leftover cell style.
function Set-TextCell($ws, $addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = 168    # Total Trades
$wsSummary.Range("B9").Value = 45.24  # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking row
# ---------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D6").Value = 65     # Trades
$wsStatus.Range("G6").Value = 47.69  # Win Rate %

# ---------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Close existing trade #168 (row 169)
$wsAll.Range("G169").Value = 0.82
$wsAll.Range("H169").Value = "CLOSED"
$wsAll.Range("K169").Value = 99.34
$wsAll.Range("L169").Value = "early_exit"
$wsAll.Range("M169").Value = 0.23

# New trade #197 (row 198) - MarketMaking, still OPEN
$wsAll.Range("A198").Value = 197
Set-TextCell $wsAll "B198" "2026-02-18"
Set-TextCell $wsAll "C198" "00:41:40"
Set-TextCell $wsAll "D198" "MarketMaking"
Set-TextCell $wsAll "E198" "DOWN"
$wsAll.Range("F198").Value = 0.82
$wsAll.Range("H198").Value = "OPEN"
$wsAll.Range("I198").Value = 0
$wsAll.Range("J198").Value = 0
$wsAll.Range("K198").Value = 99.33858346467945
$wsAll.Range("M198").Value = 0
$wsAll.Range("N198").Value = 0
$wsAll.Range("O198").Value = 0
$wsAll.Range("P198").Value = 0.6
Set-TextCell $wsAll "Q198" "Normal spread capture: 198 bps"

# New trade #198 (row 199) - EMAArbitrage, still OPEN
$wsAll.Range("A199").Value = 198
Set-TextCell $wsAll "B199" "2026-02-18"
Set-TextCell $wsAll "C199" "00:41:41"
Set-TextCell $wsAll "D199" "EMAArbitrage"
Set-TextCell $wsAll "E199" "DOWN"
$wsAll.Range("F199").Value = 0.82
$wsAll.Range("H199").Value = "OPEN"
$wsAll.Range("I199").Value = 0
$wsAll.Range("J199").Value = 0
$wsAll.Range("K199").Value = 100.270616878256
$wsAll.Range("M199").Value = 0
$wsAll.Range("N199").Value = 0
$wsAll.Range("O199").Value = 0
$wsAll.Range("P199").Value = 0.9
Set-TextCell $wsAll "Q199" "EMA:down, RSI:50.0, ROC:-45.11% | 2/3 DOWN"

# ---------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Close existing trade #168 (row 66)
$wsMM.Range("G66").Value = 0.82
$wsMM.Range("H66").Value = "CLOSED"
$wsMM.Range("K66").Value = 99.34
Set-TextCell $wsMM "P66" "early_exit"
$wsMM.Range("Q66").Value = 0.23

# New trade #197 (row 84) - still OPEN
$wsMM.Range("A84").Value = 197
Set-TextCell $wsMM "B84" "2026-02-18"
Set-TextCell $wsMM "C84" "00:41:40"
Set-TextCell $wsMM "D84" "MarketMaking"
Set-TextCell $wsMM "E84" "DOWN"
$wsMM.Range("F84").Value = 0.82
$wsMM.Range("H84").Value = "OPEN"
$wsMM.Range("I84").Value = 0
$wsMM.Range("J84").Value = 0
$wsMM.Range("K84").Value = 99.33858346467945
$wsMM.Range("L84").Value = 0
$wsMM.Range("M84").Value = 0
$wsMM.Range("N84").Value = 0.6
Set-TextCell $wsMM "O84" "Normal spread capture: 198 bps"
$wsMM.Range("Q84").Value = 0

# ---------------------------------------------------------------
# EMAArbitrage sheet
# ---------------------------------------------------------------
$wsEMA = $wb.Worksheets.Item("EMAArbitrage")

# New trade #198 (row 10) - still OPEN
$wsEMA.Range("A10").Value = 198
Set-TextCell $wsEMA "B10" "2026-02-18"
Set-TextCell $wsEMA "C10" "00:41:41"
Set-TextCell $wsEMA "D10" "EMAArbitrage"
Set-TextCell $wsEMA "E10" "DOWN"
$wsEMA.Range("F10").Value = 0.82
$wsEMA.Range("H10").Value = "OPEN"
$wsEMA.Range("I10").Value = 0
$wsEMA.Range("J10").Value = 0
$wsEMA.Range("K10").Value = 100.270616878256
$wsEMA.Range("L10").Value = 0
$wsEMA.Range("M10").Value = 0
$wsEMA.Range("N10").Value = 0.9
Set-TextCell $wsEMA "O10" "EMA:down, RSI:50.0, ROC:-45.11% | 2/3 DOWN"
$wsEMA.Range("Q10").Value = 0
